# Apply the "cryptos list" refresh (GitHub Actions update) to Sheet1.
# All Price (D) / Volume(1h) (E) columns are plain text cells, so any
# value that looks like a pure number is forced to Text via
# NumberFormat="@" before assignment (then ClearFormats() restores the
# default/general style) to avoid Excel's automatic numeric coercion
# (which would also truncate trailing zeros / change precision).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.643.39'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '3.797.73'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '707.37'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.61'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("D7").Value = '3.796.57'
$ws.Range("E7").Value = '  -1.81%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.17'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").Value = '4.438.17'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").Value = '3.726.88'
$ws.Range("E16").Value = '  -3.63%  '
$ws.Range("D17").Value = '70.690.62'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("E19").Value = '  -1.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '493.46'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.63'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.08'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.47'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").Value = '3.948.12'
$ws.Range("E28").Value = '  -1.75%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -4.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.34'
$ws.Range("D32").ClearFormats()
$ws.Range("E33").Value = '  -4.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.10'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E35").Value = '  -3.37%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.768.20'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.06'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.88%  '
$ws.Range("E39").Value = '  -3.34%  '
$ws.Range("E40").Value = '  +1.26%  '
$ws.Range("E41").Value = '  -4.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.93'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.30'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.77%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000320'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '164.40'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '422.67'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.69'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.37'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.37%  '
